$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 36
$ws.Cells.Item(3, 6).Value = 1717
$ws.Cells.Item(4, 6).Value = 1190
$ws.Cells.Item(5, 6).Value = 30
$ws.Cells.Item(7, 6).Value = 1406
$ws.Cells.Item(8, 6).Value = 69
$ws.Cells.Item(9, 6).Value = 8
$ws.Cells.Item(10, 6).Value = 104
$ws.Cells.Item(11, 6).Value = 641
$ws.Cells.Item(13, 6).Value = 110
$ws.Cells.Item(14, 6).Value = 1309
$ws.Cells.Item(15, 6).Value = 476
$ws.Cells.Item(16, 6).Value = 470
$ws.Cells.Item(19, 6).Value = 706
$ws.Cells.Item(24, 6).Value = 287
$ws.Cells.Item(27, 6).Value = 105
$ws.Cells.Item(28, 6).Value = 566
$ws.Cells.Item(31, 6).Value = 84
$ws.Cells.Item(36, 6).Value = 24

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 725
$ws.Cells.Item(15, 6).Value = 370
$ws.Cells.Item(16, 6).Value = 370
$ws.Cells.Item(26, 6).Value = 234
$ws.Cells.Item(27, 6).Value = 223

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 1752
$ws.Cells.Item(5, 6).Value = 2284
$ws.Cells.Item(6, 6).Value = 909
$ws.Cells.Item(9, 6).Value = 1132
$ws.Cells.Item(10, 6).Value = 260
$ws.Cells.Item(11, 6).Value = 77

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1752
$ws.Cells.Item(3, 6).Value = 2284
$ws.Cells.Item(4, 6).Value = 36
$ws.Cells.Item(5, 6).Value = 1717
$ws.Cells.Item(8, 6).Value = 909
$ws.Cells.Item(9, 6).Value = 1132
$ws.Cells.Item(10, 6).Value = 260
$ws.Cells.Item(11, 6).Value = 77
$ws.Cells.Item(12, 6).Value = 725
$ws.Cells.Item(13, 6).Value = 1190
$ws.Cells.Item(14, 6).Value = 30
$ws.Cells.Item(16, 6).Value = 1406
$ws.Cells.Item(18, 6).Value = 69
$ws.Cells.Item(19, 6).Value = 104
$ws.Cells.Item(20, 6).Value = 641
$ws.Cells.Item(23, 6).Value = 110
$ws.Cells.Item(24, 6).Value = 476
$ws.Cells.Item(26, 6).Value = 470
$ws.Cells.Item(28, 6).Value = 706
$ws.Cells.Item(32, 6).Value = 287
$ws.Cells.Item(35, 6).Value = 105
$ws.Cells.Item(37, 6).Value = 566
$ws.Cells.Item(39, 6).Value = 370
$ws.Cells.Item(42, 6).Value = 84
$ws.Cells.Item(46, 6).Value = 234
$ws.Cells.Item(47, 6).Value = 223

